$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Trimestre column (C) stays as plain text, not auto-converted to a date serial
$ws.Range("C2:C10").NumberFormat = "@"

# Row 2 - Santa Catarina
$ws.Range("C2").Value = "01/07/2025"
$ws.Range("D2").Value = 97.66

# Row 3 - Rondônia -> Mato Grosso
$ws.Range("A3").Value = "Mato Grosso"
$ws.Range("C3").Value = "01/07/2025"
$ws.Range("D3").Value = 97.63

# Row 4 - Mato Grosso -> Espírito Santo
$ws.Range("A4").Value = "Espírito Santo"
$ws.Range("C4").Value = "01/07/2025"
$ws.Range("D4").Value = 97.42

# Row 5 - Mato Grosso do Sul -> Rondônia
$ws.Range("A5").Value = "Rondônia"
$ws.Range("C5").Value = "01/07/2025"
$ws.Range("D5").Value = 97.38

# Row 6 - Espírito Santo -> Mato Grosso do Sul
$ws.Range("A6").Value = "Mato Grosso do Sul"
$ws.Range("C6").Value = "01/07/2025"
$ws.Range("D6").Value = 97.09999999999999

# Row 7 - Paraná
$ws.Range("C7").Value = "01/07/2025"
$ws.Range("D7").Value = 96.48

# Row 8 - Sergipe
$ws.Range("C8").Value = "01/07/2025"
$ws.Range("D8").Value = 92.26000000000001
$ws.Range("E8").Value = "23º"

# Row 9 - Brasil
$ws.Range("C9").Value = "01/07/2025"
$ws.Range("D9").Value = 94.43000000000001

# Row 10 - Nordeste
$ws.Range("C10").Value = "01/07/2025"
$ws.Range("D10").Value = 92.16
